$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.555.42"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.458.18"
$ws.Range("E4").Value = "  -1.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.92%  "
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.99%  "
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "2.838.82"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.94%  "
$ws.Range("D16").Value = "2.452.38"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.770"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("D18").Value = "41.533.31"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.87%  "
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  +5.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +4.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.09%  "
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.79%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0757"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("E39").Value = "  +5.30%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "1.958.48"
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.18%  "
$ws.Range("D48").Value = "2.698.77"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.03%  "
